$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row above row 1, shifting the existing data
# (A1:A11) down one row, to A2:A12.
$ws.Range("A1").EntireRow.Insert()

# Label the newly inserted header cell.
$ws.Range("A1").Value = "Weights"

# Mirror the recorded end-of-edit cursor position.
$ws.Range("G18").Select()
